$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new task rows to the table (E9, F10, G11)
$ws.Range("E9").Value = "Melhorar lista de abas"
$ws.Range("F10").Value = "Ao invés de Texto, grid de tabela"
$ws.Range("G11").Value = "mostrar a linha pausada"

# Update the selected/active cell to C7
$ws.Range("C7").Select()
